$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'61.815.91"
$ws.Cells.Item(2, 5).Value = "  +0.90%  "

$ws.Cells.Item(3, 4).Value = "'3.421.34"
$ws.Cells.Item(3, 5).Value = "  +3.62%  "

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "

$ws.Cells.Item(5, 4).Value = "'577.30"
$ws.Cells.Item(5, 5).Value = "  +1.98%  "

$ws.Cells.Item(6, 4).Value = "'139.60"
$ws.Cells.Item(6, 5).Value = "  +9.12%  "

$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 5).Value = "  -0.15%  "

$ws.Cells.Item(8, 4).Value = "'3.420.19"
$ws.Cells.Item(8, 5).Value = "  +3.68%  "

$ws.Cells.Item(9, 5).Value = "  +0.09%  "

$ws.Cells.Item(10, 4).Value = "'7.71"
$ws.Cells.Item(10, 5).Value = "  +5.48%  "

$ws.Cells.Item(11, 5).Value = "  +7.33%  "

$ws.Cells.Item(12, 5).Value = "  +5.86%  "

$ws.Cells.Item(13, 4).Value = "'3.995.11"
$ws.Cells.Item(13, 5).Value = "  +3.22%  "

$ws.Cells.Item(15, 4).Value = "'0.0000180"
$ws.Cells.Item(15, 5).Value = "  +7.80%  "

$ws.Cells.Item(16, 4).Value = "'3.415.14"
$ws.Cells.Item(16, 5).Value = "  +3.17%  "

$ws.Cells.Item(17, 4).Value = "'25.64"
$ws.Cells.Item(17, 5).Value = "  +6.05%  "

$ws.Cells.Item(18, 4).Value = "'61.802.96"

$ws.Cells.Item(19, 4).Value = "'14.11"
$ws.Cells.Item(19, 5).Value = "  +6.01%  "

$ws.Cells.Item(20, 4).Value = "'5.95"
$ws.Cells.Item(20, 5).Value = "  +6.01%  "

$ws.Cells.Item(21, 4).Value = "'9.49"
$ws.Cells.Item(21, 5).Value = "  +6.33%  "

$ws.Cells.Item(22, 4).Value = "'391.31"
$ws.Cells.Item(22, 5).Value = "  +10.53%  "

$ws.Cells.Item(23, 4).Value = "'0.575"
$ws.Cells.Item(23, 5).Value = "  +4.18%  "

$ws.Cells.Item(24, 4).Value = "'3.554.34"
$ws.Cells.Item(24, 5).Value = "  +3.43%  "

$ws.Cells.Item(25, 2).Value = "PEPE"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(25, 4).Value = "'0.0000127"
$ws.Cells.Item(25, 5).Value = "  +19.15%  "

$ws.Cells.Item(26, 2).Value = "Dai"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(26, 4).Value = "'1.00"
$ws.Cells.Item(26, 5).Value = "  +0.15%  "

$ws.Cells.Item(27, 4).Value = "'71.08"
$ws.Cells.Item(27, 5).Value = "  +2.97%  "

$ws.Cells.Item(28, 5).Value = "  +15.30%  "

$ws.Cells.Item(29, 4).Value = "'7.83"
$ws.Cells.Item(29, 5).Value = "  +9.98%  "

$ws.Cells.Item(30, 4).Value = "'0.992"
$ws.Cells.Item(30, 5).Value = "  -0.70%  "

$ws.Cells.Item(31, 4).Value = "'8.33"
$ws.Cells.Item(31, 5).Value = "  +7.37%  "

$ws.Cells.Item(32, 4).Value = "'0.159"
$ws.Cells.Item(32, 5).Value = "  +7.05%  "

$ws.Cells.Item(33, 5).Value = "  +2.60%  "

$ws.Cells.Item(34, 4).Value = "'3.449.75"
$ws.Cells.Item(34, 5).Value = "  +3.56%  "

$ws.Cells.Item(35, 5).Value = "  -0.10%  "

$ws.Cells.Item(36, 4).Value = "'23.72"
$ws.Cells.Item(36, 5).Value = "  +5.04%  "

$ws.Cells.Item(37, 4).Value = "'5.54"
$ws.Cells.Item(37, 5).Value = "  +5.69%  "

$ws.Cells.Item(38, 4).Value = "'7.07"
$ws.Cells.Item(38, 5).Value = "  +4.76%  "

$ws.Cells.Item(39, 5).Value = "  +6.67%  "

$ws.Cells.Item(40, 4).Value = "'162.14"
$ws.Cells.Item(40, 5).Value = "  -0.40%  "

$ws.Cells.Item(41, 4).Value = "'0.0799"
$ws.Cells.Item(41, 5).Value = "  +6.43%  "

$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(42, 4).Value = "'1.74"
$ws.Cells.Item(42, 5).Value = "  +12.13%  "

$ws.Cells.Item(43, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(43, 4).Value = "'1.00"
$ws.Cells.Item(43, 5).Value = "  +0.03%  "

$ws.Cells.Item(44, 2).Value = "Filecoin"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(44, 4).Value = "'4.50"
$ws.Cells.Item(44, 5).Value = "  +3.04%  "

$ws.Cells.Item(45, 2).Value = "Mantle"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(45, 4).Value = "'0.778"
$ws.Cells.Item(45, 5).Value = "  +5.05%  "

$ws.Cells.Item(46, 2).Value = "ONDO"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(46, 4).Value = "'1.23"
$ws.Cells.Item(46, 5).Value = "  +9.71%  "

$ws.Cells.Item(47, 4).Value = "'41.27"
$ws.Cells.Item(47, 5).Value = "  +0.42%  "

$ws.Cells.Item(48, 4).Value = "'23.54"
$ws.Cells.Item(48, 5).Value = "  +6.05%  "

$ws.Cells.Item(49, 4).Value = "'7.04"
$ws.Cells.Item(49, 5).Value = "  +5.78%  "

$ws.Cells.Item(50, 4).Value = "'23.00"
$ws.Cells.Item(50, 5).Value = "  +8.99%  "

$ws.Cells.Item(51, 4).Value = "'2.360.66"
$ws.Cells.Item(51, 5).Value = "  +9.83%  "
